$d = $word.ActiveDocument

function Get-ParagraphIndexStartingWith($doc, $prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return $null
}

# --- Hunk 1 -----------------------------------------------------------------
# Right after the paragraph "Aggregation renders Resources into Functional
# Template Form ..." sits an empty paragraph (pBdr/shd formatting, no text).
# Two new paragraphs are inserted right after that empty paragraph:
#   1) another empty paragraph, same formatting, but with ind left=0/firstLine=0
#   2) a paragraph with the same ind, holding the new sentence.

$aggIdx = Get-ParagraphIndexStartingWith $d "Aggregation renders Resources into Functional Template Form"
$emptyAfterAgg = $d.Paragraphs.Item($aggIdx + 1)

$emptyAfterAgg.Range.InsertParagraphAfter()
$hunk1Blank = $d.Paragraphs.Item($aggIdx + 2)
$hunk1Blank.Format.LeftIndent = 0
$hunk1Blank.Format.FirstLineIndent = 0

$hunk1Blank.Range.InsertParagraphAfter()
$hunk1Text = $d.Paragraphs.Item($aggIdx + 3)
$hunk1Text.Range.Text = "Resources aggregate into Kinds. Kinds aggregate into Statements, Statements aggregate into Mappings. Mappings aggregate into Transforms. Hierarchy aligns Wrapper types."
$hunk1Text.Format.LeftIndent = 0
$hunk1Text.Format.FirstLineIndent = 0

# --- Hunk 2 -------------------------------------------------------------------
# Fix the missing close-paren in "Inputs: (Class, Instance, Attribute, Value"
# then add two new numbered-list paragraphs right after that paragraph:
#   1) a paragraph with the same new sentence as above
#   2) a trailing empty numbered-list paragraph

$d.Content.Find.Execute(
    "Inputs: (Class, Instance, Attribute, Value: Aggregate: Functional Form implemented in Quad interface.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Inputs: (Class, Instance, Attribute, Value): Aggregate: Functional Form implemented in Quad interface.",
    2)

$inputsIdx = Get-ParagraphIndexStartingWith $d "Inputs: (Class, Instance, Attribute, Value): Aggregate"
$inputsP = $d.Paragraphs.Item($inputsIdx)

$inputsP.Range.InsertParagraphAfter()
$hunk2Text = $d.Paragraphs.Item($inputsIdx + 1)
$hunk2Text.Range.Text = "Resources aggregate into Kinds. Kinds aggregate into Statements, Statements aggregate into Mappings. Mappings aggregate into Transforms. Hierarchy aligns Wrapper types."

$hunk2Text.Range.InsertParagraphAfter()
